$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dynamic user-id based parameters are replaced with a placeholder that
# is resolved from system properties at test run time, instead of being a
# hard-coded user id.
$oldHeader = "X-1P-User=a99ba4dc-45be-4ad2-9c9e-22e78584b82b||Content-Type=application/json"
$newHeader = "X-1P-User=(SYS_USER1)||Content-Type=application/json"

$oldPath = "/comments/user/a99ba4dc-45be-4ad2-9c9e-22e78584b82b"
$newPath = "/comments/user/(SYS_USER1)"

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq $oldHeader) {
            $cell.Value = $newHeader
        } elseif ($val -eq $oldPath) {
            $cell.Value = $newPath
        }
    }
}

# Tidy up the left-over blank placeholder cells (no longer needed) so they
# don't linger as empty cell records.
$blankCells = @("G3","I3","K3","G4","K4","G5","K5","G6","K6","G7","K7")
foreach ($addr in $blankCells) {
    $ws.Range($addr).ClearContents()
}

# Update the selection / active cell shown when the workbook is reopened.
$ws.Range("A7").Select()
